$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 31000
$ws.Range("I81").Value = 20000
$ws.Range("J81").Value = 36500
$ws.Range("K81").Value = 20000
$ws.Range("L81").Value = 36500
$ws.Range("M81").Value = -19002
$ws.Range("N81").Value = -38496
$ws.Range("H84").Value = 31000
$ws.Range("I84").Value = 20000
$ws.Range("J84").Value = 36500
$ws.Range("K84").Value = 60000
$ws.Range("L84").Value = 109500
$ws.Range("M84").Value = -55008
$ws.Range("N84").Value = -119484
$ws.Range("H113").Value = 107422400
$ws.Range("I113").Value = 27779882
$ws.Range("J113").Value = 136383330
$ws.Range("K113").Value = 27779882
$ws.Range("L113").Value = 136383330
$ws.Range("M113").Value = -27776628
$ws.Range("N113").Value = -136389838
$ws.Range("H137").Value = 5824.516
$ws.Range("I137").Value = 2714.7778
$ws.Range("K137").Value = 8144.3334
$ws.Range("M137").Value = -5594.3334
$ws.Range("H138").Value = 1102041
$ws.Range("J138").Value = 2444463
$ws.Range("L138").Value = 7333389
$ws.Range("N138").Value = -7343669

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2784051
$ws.Range("I32").Value = 2946178.5
$ws.Range("K32").Value = 2946178.5
$ws.Range("M32").Value = -2945891.5
$ws.Range("H45").Value = 4636.8696
$ws.Range("I45").Value = 1842.8889
$ws.Range("K45").Value = 1842.8889
$ws.Range("M45").Value = -1465.8889
$ws.Range("H61").Value = 41673972
$ws.Range("I61").Value = 3459.9375
$ws.Range("K61").Value = 3459.9375
$ws.Range("M61").Value = -3247.9375
$ws.Range("H74").Value = 40801.703
$ws.Range("I74").Value = 69316.92999999999
$ws.Range("K74").Value = 69316.92999999999
$ws.Range("M74").Value = -68442.92999999999
$ws.Range("H77").Value = 40801.703
$ws.Range("I77").Value = 69316.92999999999
$ws.Range("K77").Value = 346584.65
$ws.Range("M77").Value = -342216.65
$ws.Range("H95").Value = 38562.6
$ws.Range("J95").Value = 38562.6
$ws.Range("L95").Value = 38562.6
$ws.Range("N95").Value = -44054.6
$ws.Range("H110").Value = 27779010
$ws.Range("I110").Value = 1217.8889
$ws.Range("K110").Value = 1217.8889
$ws.Range("M110").Value = 827.1111000000001
$ws.Range("H132").Value = 4599.1133
$ws.Range("I132").Value = 2189.484
$ws.Range("K132").Value = 6568.451999999999
$ws.Range("M132").Value = -4038.451999999999
$ws.Range("H136").Value = 41673972
$ws.Range("I136").Value = 3459.9375
$ws.Range("K136").Value = 10379.8125
$ws.Range("M136").Value = -7829.8125
$ws.Range("H139").Value = 68715
$ws.Range("J139").Value = 68715
$ws.Range("L139").Value = 68715
$ws.Range("N139").Value = -78995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 11913119
$ws.Range("I20").Value = 18525830
$ws.Range("J20").Value = 10241.6
$ws.Range("K20").Value = 18525830
$ws.Range("L20").Value = 10241.6
$ws.Range("M20").Value = -18525583
$ws.Range("N20").Value = -10735.6
$ws.Range("H128").Value = 4629.25
$ws.Range("I128").Value = 4629.25
$ws.Range("K128").Value = 13887.75
$ws.Range("M128").Value = -11397.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8146.234
$ws.Range("I31").Value = 3628.182
$ws.Range("J31").Value = 9526.75
$ws.Range("K31").Value = 3628.182
$ws.Range("L31").Value = 9526.75
$ws.Range("M31").Value = -3333.182
$ws.Range("N31").Value = -10116.75
$ws.Range("H34").Value = 8146.234
$ws.Range("I34").Value = 3628.182
$ws.Range("J34").Value = 9526.75
$ws.Range("K34").Value = 3628.182
$ws.Range("L34").Value = 9526.75
$ws.Range("M34").Value = -3426.182
$ws.Range("N34").Value = -9930.75
$ws.Range("H76").Value = 4718.5
$ws.Range("I76").Value = 4718.5
$ws.Range("K76").Value = 4718.5
$ws.Range("M76").Value = -4403.5
$ws.Range("H79").Value = 4718.5
$ws.Range("I79").Value = 4718.5
$ws.Range("K79").Value = 4718.5
$ws.Range("M79").Value = -3626.5
$ws.Range("H86").Value = 12024556
$ws.Range("I86").Value = 22327782
$ws.Range("K86").Value = 22327782
$ws.Range("M86").Value = -22326659
$ws.Range("H89").Value = 12024556
$ws.Range("I89").Value = 22327782
$ws.Range("K89").Value = 111638910
$ws.Range("M89").Value = -111633294
$ws.Range("H96").Value = 9078.333000000001
$ws.Range("J96").Value = 9078.333000000001
$ws.Range("L96").Value = 9078.333000000001
$ws.Range("N96").Value = -14570.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1607.625
$ws.Range("I5").Value = 724.25
$ws.Range("K5").Value = 2172.75
$ws.Range("M5").Value = -2060.75
$ws.Range("H56").Value = 7125
$ws.Range("I56").Value = 7125
$ws.Range("K56").Value = 7125
$ws.Range("M56").Value = -6595
$ws.Range("H113").Value = 5858.5
$ws.Range("J113").Value = 9280
$ws.Range("L113").Value = 27840
$ws.Range("N113").Value = -32180
$ws.Range("H131").Value = 45020.22
$ws.Range("I131").Value = 998.75
$ws.Range("J131").Value = 54287.895
$ws.Range("K131").Value = 2996.25
$ws.Range("L131").Value = 162863.685
$ws.Range("M131").Value = 2043.75
$ws.Range("N131").Value = -172943.685
$ws.Range("H134").Value = 46438.668
$ws.Range("I134").Value = 54373.65
$ws.Range("K134").Value = 163120.95
$ws.Range("M134").Value = -158050.95
$ws.Range("H135").Value = 1607.625
$ws.Range("I135").Value = 724.25
$ws.Range("K135").Value = 6518.25
$ws.Range("M135").Value = -3983.25
$ws.Range("H137").Value = 69511.97
$ws.Range("J137").Value = 65745.875
$ws.Range("L137").Value = 197237.625
$ws.Range("N137").Value = -207437.625
$ws.Range("H138").Value = 85073.62
$ws.Range("J138").Value = 12166.667
$ws.Range("L138").Value = 36500.001
$ws.Range("N138").Value = -46780.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3206.739
$ws.Range("I102").Value = 2893.0952
$ws.Range("K102").Value = 2893.0952
$ws.Range("M102").Value = -1271.0952
$ws.Range("H113").Value = 2713.4666
$ws.Range("I113").Value = 2785.4285
$ws.Range("J113").Value = 2650.5
$ws.Range("K113").Value = 2785.4285
$ws.Range("L113").Value = 2650.5
$ws.Range("M113").Value = -615.4285
$ws.Range("N113").Value = -6990.5
$ws.Range("H122").Value = 5573771.5
$ws.Range("J122").Value = 998
$ws.Range("L122").Value = 2994
$ws.Range("N122").Value = -7894
$ws.Range("H126").Value = 4920
$ws.Range("I126").Value = 2161.2727
$ws.Range("J126").Value = 7448.8335
$ws.Range("K126").Value = 6483.8181
$ws.Range("L126").Value = 22346.5005
$ws.Range("M126").Value = -4013.8181
$ws.Range("N126").Value = -27286.5005
$ws.Range("H132").Value = 10002.375
$ws.Range("I132").Value = 3008
$ws.Range("J132").Value = 14199
$ws.Range("K132").Value = 9024
$ws.Range("L132").Value = 42597
$ws.Range("M132").Value = -6494
$ws.Range("N132").Value = -47657

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 2550000
$ws.Range("I20").Value = 2600000
$ws.Range("K20").Value = 2600000
$ws.Range("M20").Value = -2599774
$ws.Range("H61").Value = 3488.3704
$ws.Range("I61").Value = 1893.2858
$ws.Range("J61").Value = 5206.154
$ws.Range("K61").Value = 1893.2858
$ws.Range("L61").Value = 5206.154
$ws.Range("M61").Value = -1691.2858
$ws.Range("N61").Value = -5610.154
$ws.Range("H113").Value = 3488.3704
$ws.Range("I113").Value = 1893.2858
$ws.Range("J113").Value = 5206.154
$ws.Range("K113").Value = 1893.2858
$ws.Range("L113").Value = 5206.154
$ws.Range("M113").Value = 276.7141999999999
$ws.Range("N113").Value = -9546.154
$ws.Range("H132").Value = 12202273
$ws.Range("I132").Value = 27781556
$ws.Range("K132").Value = 83344668
$ws.Range("M132").Value = -83342138
$ws.Range("H136").Value = 12659.444
$ws.Range("I136").Value = 4018.5625
$ws.Range("J136").Value = 17426.828
$ws.Range("K136").Value = 12055.6875
$ws.Range("L136").Value = 52280.484
$ws.Range("M136").Value = -9505.6875
$ws.Range("N136").Value = -57380.484

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 25876.5
$ws.Range("J69").Value = 28106
$ws.Range("L69").Value = 28106
$ws.Range("N69").Value = -29604
$ws.Range("H72").Value = 25876.5
$ws.Range("J72").Value = 28106
$ws.Range("L72").Value = 84318
$ws.Range("N72").Value = -91806
$ws.Range("H80").Value = 22109.715
$ws.Range("I80").Value = 18591
$ws.Range("J80").Value = 24748.75
$ws.Range("K80").Value = 18591
$ws.Range("L80").Value = 24748.75
$ws.Range("M80").Value = -17593
$ws.Range("N80").Value = -26744.75
$ws.Range("H83").Value = 22109.715
$ws.Range("I83").Value = 18591
$ws.Range("J83").Value = 24748.75
$ws.Range("K83").Value = 55773
$ws.Range("L83").Value = 74246.25
$ws.Range("M83").Value = -50781
$ws.Range("N83").Value = -84230.25
$ws.Range("H132").Value = 5353.512
$ws.Range("I132").Value = 5113.3447
$ws.Range("K132").Value = 15340.0341
$ws.Range("M132").Value = -12810.0341
$ws.Range("H136").Value = 21493442
$ws.Range("I136").Value = 50002376
$ws.Range("K136").Value = 150007128
$ws.Range("M136").Value = -150004578
